# Update countries & provincias Spain
# Refresh the COVID case counters for a handful of countries, update the
# "last updated" timestamp caption, then re-sort the data range by total
# cases (column B) descending - which is how the sheet is always kept - so
# any country whose count jumped past its neighbours naturally moves to its
# new rank (e.g. "Islas Turcas y Caicos" overtakes "San Martin (Parte
# Holandesa)" and "Papua Nueva Guinea").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($Country, $TotalCases, $NewCases, $Active, $Recovered, $Critical, $Deaths) {
    $found = $ws.Range("A4:A219").Find($Country)
    $r = $found.Row

    $ws.Cells.Item($r, 2).Value = $TotalCases
    $ws.Cells.Item($r, 3).Value = $NewCases
    if ($Active -ne $null) { $ws.Cells.Item($r, 4).Value = $Active }
    $ws.Cells.Item($r, 5).Value = $Recovered
    if ($Critical -ne $null) { $ws.Cells.Item($r, 7).Value = $Critical }
    if ($Deaths -ne $null) { $ws.Cells.Item($r, 8).Value = $Deaths }
}

Set-CountryRow "Bolivia"    108427 992 43887 60098 76   4442
Set-CountryRow "Kazajistan" 104543 230 $null  13416 $null $null
Set-CountryRow "Belgica"    81468  574 18204  53276 3    9988
Set-CountryRow "Honduras"   53983  602 8449   43891 11   1643
Set-CountryRow "Uzbekistan" 38532  120 34576  3689  2    267
Set-CountryRow "Australia"  24811  209 $null  5550  17   502
Set-CountryRow "Islas Turcas y Caicos" 383 36 102 279 $null 2
Set-CountryRow "Fiyi"       28     0   23     4     $null $null

# Update the "last updated" caption.
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 05:23"

# Re-sort the table (rows 4:219) by total cases (column B) descending, as
# the sheet is normally kept, so rank changes ripple through automatically.
$dataRange = $ws.Range("A4:H219")
$keyRange = $ws.Range("B4:B219")
$dataRange.Sort($keyRange, 2, $null, $null, 1, $null, 1, 1)
